$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.043.15"
$ws.Range("E2").Value = "  -1.66%  "

$ws.Range("D3").Value = "2.304.27"
$ws.Range("E3").Value = "  -1.95%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.90%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.69%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.608"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.71%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0905"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.11%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.49"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.90%  "

$ws.Range("E13").Value = "  +0.91%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.973"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.70%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.54%  "

$ws.Range("D16").Value = "2.653.29"
$ws.Range("E16").Value = "  -2.07%  "

$ws.Range("D17").Value = "2.309.90"
$ws.Range("E17").Value = "  -1.35%  "

$ws.Range("D18").Value = "42.080.61"
$ws.Range("E18").Value = "  -1.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.29%  "

$ws.Range("E20").Value = "  -0.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.14%  "

$ws.Range("E22").Value = "  -0.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "268.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.25%  "

$ws.Range("E24").Value = "  -0.81%  "

$ws.Range("E25").Value = "  +4.53%  "

$ws.Range("E26").Value = "  +0.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.79%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.27%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "23.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "165.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.53%  "

$ws.Range("E32").Value = "  -0.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.11%  "

$ws.Range("E35").Value = "  +1.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0353"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.60%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "103.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +17.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.47%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.227"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.39%  "

$ws.Range("E45").Value = "  +0.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "115.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "78.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.95%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.09%  "

$ws.Range("E51").Value = "  +2.30%  "
